$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted as row 453 ("Región de Coquimbo",
# 2023-01-05), pushing the existing rows 453-485 down to 454-486 and
# carrying the row directly below row 485 (the former last record) into
# the newly created row 486.
$ws.Rows.Item(453).Insert()

$ws.Range("A453").Value = 4
$ws.Range("B453").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C453").Value = "Los Lagos"
$ws.Range("D453").Value = 44931
$ws.Range("E453").Value = 10
$ws.Range("F453").Value = 100114013
$ws.Range("G453").Value = "Zanahoria"
$ws.Range("H453").Value = "Sin especificar"
$ws.Range("I453").Value = "Primera"
$ws.Range("J453").Value = 300
$ws.Range("K453").Value = 15000
$ws.Range("L453").Value = 16000
$ws.Range("M453").Value = 15500
$ws.Range("N453").Value = "`$/saco 20 kilos"
$ws.Range("O453").Value = "Región de Coquimbo"
$ws.Range("P453").Value = 775
$ws.Range("Q453").Value = 20
$ws.Range("R453").Value = "Hortaliza"
